# Apply the "adding XR UI" edit to buttonConfiguration.xlsx
# - fills in a few previously-empty cells on the "plan" sheet with new
#   button-mapping labels
# - renames "gripper control?" -> "gripper control" in its cell
# - moves the active selection on the "plan" sheet from G4 to E6

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("plan")

# New short-key / mode labels added to previously blank cells
$ws.Range("F3").Value  = "select"
$ws.Range("E7").Value  = "align short key"
$ws.Range("E8").Value  = "switch mode"
$ws.Range("F8").Value  = "switch mode"

# "gripper control?" -> "gripper control"
$ws.Range("E12").Value = "gripper control"

# Update the active cell/selection on the bottom-right frozen pane
$ws.Range("E6").Select()
